$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3246.077
$ws.Range("I15").Value = 3246.077
$ws.Range("K15").Value = 9738.231
$ws.Range("M15").Value = -9569.231

$ws.Range("H106").Value = 1159.5625
$ws.Range("I106").Value = 1118.0714
$ws.Range("J106").Value = 1450
$ws.Range("K106").Value = 1118.0714
$ws.Range("L106").Value = 1450
$ws.Range("M106").Value = -487.0714
$ws.Range("N106").Value = -2712

$ws.Range("H107").Value = 29854.176
$ws.Range("I107").Value = 38834.77
$ws.Range("J107").Value = 667.25
$ws.Range("K107").Value = 38834.77
$ws.Range("L107").Value = 667.25
$ws.Range("M107").Value = -36914.77
$ws.Range("N107").Value = -4507.25

$ws.Range("H137").Value = 1313.4897
$ws.Range("I137").Value = 1028.2122
$ws.Range("J137").Value = 1901.875
$ws.Range("K137").Value = 3084.6366
$ws.Range("L137").Value = 5705.625
$ws.Range("M137").Value = -534.6365999999998
$ws.Range("N137").Value = -10805.625

$ws.Range("H138").Value = 4061.66
$ws.Range("I138").Value = 1197.7273
$ws.Range("J138").Value = 9621.058999999999
$ws.Range("K138").Value = 3593.1819
$ws.Range("L138").Value = 28863.177
$ws.Range("M138").Value = 1546.8181
$ws.Range("N138").Value = -39143.177

$ws.Range("H141").Value = 3773.6416
$ws.Range("I141").Value = 1983.5
$ws.Range("J141").Value = 20959
$ws.Range("K141").Value = 5950.5
$ws.Range("L141").Value = 62877
$ws.Range("M141").Value = -770.5
$ws.Range("N141").Value = -73237

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10841.077
$ws.Range("I32").Value = 11646.302
$ws.Range("K32").Value = 11646.302
$ws.Range("M32").Value = -11359.302

$ws.Range("H61").Value = 1259
$ws.Range("I61").Value = 1091.5862
$ws.Range("J61").Value = 1700.3636
$ws.Range("K61").Value = 1091.5862
$ws.Range("L61").Value = 1700.3636
$ws.Range("M61").Value = -879.5862
$ws.Range("N61").Value = -2124.3636

$ws.Range("H136").Value = 1259
$ws.Range("I136").Value = 1091.5862
$ws.Range("J136").Value = 1700.3636
$ws.Range("K136").Value = 3274.7586
$ws.Range("L136").Value = 5101.0908
$ws.Range("M136").Value = -724.7586000000001
$ws.Range("N136").Value = -10201.0908

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1460.2759
$ws.Range("I134").Value = 1336.9814
$ws.Range("J134").Value = 3124.75
$ws.Range("K134").Value = 4010.9442
$ws.Range("L134").Value = 9374.25
$ws.Range("M134").Value = -1475.9442
$ws.Range("N134").Value = -14444.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1454.3
$ws.Range("I31").Value = 1260.8206
$ws.Range("J31").Value = 9000
$ws.Range("K31").Value = 1260.8206
$ws.Range("L31").Value = 9000
$ws.Range("M31").Value = -965.8206
$ws.Range("N31").Value = -9590

$ws.Range("H34").Value = 1454.3
$ws.Range("I34").Value = 1260.8206
$ws.Range("J34").Value = 9000
$ws.Range("K34").Value = 1260.8206
$ws.Range("L34").Value = 9000
$ws.Range("M34").Value = -1058.8206
$ws.Range("N34").Value = -9404

$ws.Range("H132").Value = 222652.33
$ws.Range("I132").Value = 288354.1
$ws.Range("J132").Value = 2082.0715
$ws.Range("K132").Value = 865062.2999999999
$ws.Range("L132").Value = 6246.2145
$ws.Range("M132").Value = -862532.2999999999
$ws.Range("N132").Value = -11306.2145

$ws.Range("H134").Value = 1255.5493
$ws.Range("I134").Value = 1132.862
$ws.Range("J134").Value = 1802.9231
$ws.Range("K134").Value = 3398.586
$ws.Range("L134").Value = 5408.7693
$ws.Range("M134").Value = -863.5860000000002
$ws.Range("N134").Value = -10478.7693

$ws.Range("H135").Value = 48000
$ws.Range("J135").Value = 48000
$ws.Range("L135").Value = 48000
$ws.Range("N135").Value = -58140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 55562116
$ws.Range("I3").Value = 100003610
$ws.Range("J3").Value = 10254.125
$ws.Range("K3").Value = 300010830
$ws.Range("L3").Value = 30762.375
$ws.Range("M3").Value = -300010718
$ws.Range("N3").Value = -30986.375

$ws.Range("H5").Value = 2254.1428
$ws.Range("I5").Value = 2560.8
$ws.Range("K5").Value = 7682.400000000001
$ws.Range("M5").Value = -7570.400000000001

$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()

$ws.Range("H23").Value = 427.74075
$ws.Range("I23").Value = 293.33334
$ws.Range("J23").Value = 444.54166
$ws.Range("K23").Value = 880.0000200000001
$ws.Range("L23").Value = 1333.62498
$ws.Range("M23").Value = -645.0000200000001
$ws.Range("N23").Value = -1803.62498

$ws.Range("H25").Value = 2674.75
$ws.Range("I25").Value = 100
$ws.Range("J25").Value = 3533
$ws.Range("K25").Value = 300
$ws.Range("L25").Value = 10599
$ws.Range("M25").Value = -131
$ws.Range("N25").Value = -10937

$ws.Range("H30").Value = 2674.75
$ws.Range("I30").Value = 100
$ws.Range("J30").Value = 3533
$ws.Range("K30").Value = 300
$ws.Range("L30").Value = 10599
$ws.Range("M30").Value = -198
$ws.Range("N30").Value = -10803

$ws.Range("H34").Value = 522.1429000000001
$ws.Range("I34").Value = 200
$ws.Range("J34").Value = 592.1739
$ws.Range("K34").Value = 600
$ws.Range("L34").Value = 1776.5217
$ws.Range("M34").Value = -516
$ws.Range("N34").Value = -1944.5217

$ws.Range("H39").Value = 2753.3333
$ws.Range("J39").Value = 2753.3333
$ws.Range("L39").Value = 8259.999899999999
$ws.Range("N39").Value = -8847.999899999999

$ws.Range("H55").Value = 5449.1
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 5449.1
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 16347.3
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -16701.3

$ws.Range("H135").Value = 2254.1428
$ws.Range("I135").Value = 2560.8
$ws.Range("K135").Value = 23047.2
$ws.Range("M135").Value = -20512.2

$ws.Range("H139").Value = 2129.5
$ws.Range("I139").Value = 1826.3158
$ws.Range("J139").Value = 2769.5557
$ws.Range("K139").Value = 5478.9474
$ws.Range("L139").Value = 8308.667099999999
$ws.Range("M139").Value = -338.9474
$ws.Range("N139").Value = -18588.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2618.7183
$ws.Range("I132").Value = 2161.3447
$ws.Range("J132").Value = 4659.3076
$ws.Range("K132").Value = 6484.034100000001
$ws.Range("L132").Value = 13977.9228
$ws.Range("M132").Value = -3954.034100000001
$ws.Range("N132").Value = -19037.9228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1950
$ws.Range("I96").Value = 2333.6667
$ws.Range("J96").Value = 799
$ws.Range("K96").Value = 2333.6667
$ws.Range("L96").Value = 799
$ws.Range("M96").Value = -960.6667000000002
$ws.Range("N96").Value = -3545

$ws.Range("H132").Value = 1024.3704
$ws.Range("I132").Value = 739.8333
$ws.Range("J132").Value = 1593.4445
$ws.Range("K132").Value = 2219.4999
$ws.Range("L132").Value = 4780.333500000001
$ws.Range("M132").Value = 310.5001000000002
$ws.Range("N132").Value = -9840.333500000001

$ws.Range("H136").Value = 2140.2778
$ws.Range("I136").Value = 1725.862
$ws.Range("J136").Value = 3857.1428
$ws.Range("K136").Value = 5177.586
$ws.Range("L136").Value = 11571.4284
$ws.Range("M136").Value = -2627.586
$ws.Range("N136").Value = -16671.4284
